$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5153.773
$ws.Range("I19").Value = 8237.923000000001
$ws.Range("J19").Value = 698.8889
$ws.Range("K19").Value = 8237.923000000001
$ws.Range("L19").Value = 698.8889
$ws.Range("M19").Value = -8062.923000000001
$ws.Range("N19").Value = -1048.8889

$ws.Range("H33").Value = 4973.091
$ws.Range("I33").Value = 6838.3335
$ws.Range("J33").Value = 976.1429000000001
$ws.Range("K33").Value = 6838.3335
$ws.Range("L33").Value = 976.1429000000001
$ws.Range("M33").Value = -6609.3335
$ws.Range("N33").Value = -1434.1429

$ws.Range("H129").Value = 23810352
$ws.Range("I129").Value = 166666960
$ws.Range("J129").Value = 916.19446
$ws.Range("K129").Value = 500000880
$ws.Range("L129").Value = 2748.58338
$ws.Range("M129").Value = -499995880
$ws.Range("N129").Value = -12748.58338

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 29178
$ws.Range("J7").Value = 29178
$ws.Range("L7").Value = 29178
$ws.Range("N7").Value = -29406

$ws.Range("H32").Value = 4831.904
$ws.Range("I32").Value = 3005.2856
$ws.Range("J32").Value = 34666.668
$ws.Range("K32").Value = 3005.2856
$ws.Range("L32").Value = 34666.668
$ws.Range("M32").Value = -2718.2856
$ws.Range("N32").Value = -35240.668

$ws.Range("H44").Value = 41688.89
$ws.Range("J44").Value = 41688.89
$ws.Range("L44").Value = 41688.89
$ws.Range("N44").Value = -42664.89

$ws.Range("H62").Value = 31666.666
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 31666.666
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 31666.666
$ws.Range("N62").Value = -32914.666
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 31666.666
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 31666.666
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 94999.99800000001
$ws.Range("N65").Value = -101239.998
$ws.Range("M65").ClearContents()

$ws.Range("H68").Value = 63155
$ws.Range("J68").Value = 63155
$ws.Range("L68").Value = 63155
$ws.Range("N68").Value = -64777

$ws.Range("H71").Value = 63155
$ws.Range("J71").Value = 63155
$ws.Range("L71").Value = 189465
$ws.Range("N71").Value = -197577

$ws.Range("H81").Value = 90000
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H82").Value = 30000
$ws.Range("J82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30722

$ws.Range("H84").Value = 90000
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H85").Value = 30000
$ws.Range("J85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32496

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H87").Value = 40000
$ws.Range("J87").Value = 40000
$ws.Range("L87").Value = 40000
$ws.Range("N87").Value = -42496

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H90").Value = 40000
$ws.Range("J90").Value = 40000
$ws.Range("L90").Value = 120000
$ws.Range("N90").Value = -132480

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2146.5186
$ws.Range("I20").Value = 2483.5
$ws.Range("J20").Value = 1656.3636
$ws.Range("K20").Value = 2483.5
$ws.Range("L20").Value = 1656.3636
$ws.Range("M20").Value = -2236.5
$ws.Range("N20").Value = -2150.3636

$ws.Range("H107").Value = 851.2
$ws.Range("I107").Value = 812.44446
$ws.Range("K107").Value = 812.44446
$ws.Range("M107").Value = 1107.55554

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 943.4375
$ws.Range("I22").Value = 1109.3846
$ws.Range("J22").Value = 224.33333
$ws.Range("K22").Value = 1109.3846
$ws.Range("L22").Value = 224.33333
$ws.Range("M22").Value = -759.3846000000001
$ws.Range("N22").Value = -924.3333299999999

$ws.Range("H31").Value = 14875.12
$ws.Range("I31").Value = 1331.7858
$ws.Range("J31").Value = 32112.092
$ws.Range("K31").Value = 1331.7858
$ws.Range("L31").Value = 32112.092
$ws.Range("M31").Value = -1036.7858
$ws.Range("N31").Value = -32702.092

$ws.Range("H34").Value = 14875.12
$ws.Range("I34").Value = 1331.7858
$ws.Range("J34").Value = 32112.092
$ws.Range("K34").Value = 1331.7858
$ws.Range("L34").Value = 32112.092
$ws.Range("M34").Value = -1129.7858
$ws.Range("N34").Value = -32516.092

$ws.Range("H94").Value = 3384.8572
$ws.Range("I94").Value = 3963.3333
$ws.Range("J94").Value = 2951
$ws.Range("K94").Value = 3963.3333
$ws.Range("L94").Value = 2951
$ws.Range("M94").Value = -3512.3333
$ws.Range("N94").Value = -3853

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1229.8
$ws.Range("I5").Value = 478.73685
$ws.Range("K5").Value = 1436.21055
$ws.Range("M5").Value = -1324.21055

$ws.Range("H12").Value = 85.8125
$ws.Range("I12").Value = 3.142857
$ws.Range("K12").Value = 9.428571
$ws.Range("M12").Value = 163.571429

$ws.Range("H96").Value = 6666
$ws.Range("J96").Value = 6666
$ws.Range("L96").Value = 19998
$ws.Range("N96").Value = -24116

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H113").Value = 1168.6571
$ws.Range("I113").Value = 1373.7273
$ws.Range("J113").Value = 1074.6666
$ws.Range("K113").Value = 4121.1819
$ws.Range("L113").Value = 3223.9998
$ws.Range("M113").Value = -1951.1819
$ws.Range("N113").Value = -7563.9998

$ws.Range("H122").Value = 783.125
$ws.Range("I122").Value = 641.4545000000001
$ws.Range("J122").Value = 1094.8
$ws.Range("K122").Value = 5773.0905
$ws.Range("L122").Value = 9853.199999999999
$ws.Range("M122").Value = -3323.0905
$ws.Range("N122").Value = -14753.2

$ws.Range("H131").Value = 1341.38
$ws.Range("I131").Value = 640
$ws.Range("K131").Value = 1920
$ws.Range("M131").Value = 3120

$ws.Range("H135").Value = 1229.8
$ws.Range("I135").Value = 478.73685
$ws.Range("K135").Value = 4308.63165
$ws.Range("M135").Value = -1773.63165

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 18000
$ws.Range("J43").Value = 18000
$ws.Range("L43").Value = 18000
$ws.Range("N43").Value = -18302

$ws.Range("H70").Value = 6836.0586
$ws.Range("I70").Value = 8392.182000000001
$ws.Range("J70").Value = 3983.1667
$ws.Range("K70").Value = 8392.182000000001
$ws.Range("L70").Value = 3983.1667
$ws.Range("M70").Value = -8122.182000000001
$ws.Range("N70").Value = -4523.1667

$ws.Range("H73").Value = 6836.0586
$ws.Range("I73").Value = 8392.182000000001
$ws.Range("J73").Value = 3983.1667
$ws.Range("K73").Value = 8392.182000000001
$ws.Range("L73").Value = 3983.1667
$ws.Range("M73").Value = -7456.182000000001
$ws.Range("N73").Value = -5855.1667

$ws.Range("H113").Value = 1526.1428
$ws.Range("I113").Value = 1592.75
$ws.Range("J113").Value = 1437.3334
$ws.Range("K113").Value = 1592.75
$ws.Range("L113").Value = 1437.3334
$ws.Range("M113").Value = 577.25
$ws.Range("N113").Value = -5777.3334

$ws.Range("H116").Value = 45000
$ws.Range("J116").Value = 45000
$ws.Range("L116").Value = 45000
$ws.Range("N116").Value = -54178

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1055.963
$ws.Range("I22").Value = 855.7222
$ws.Range("J22").Value = 1456.4445
$ws.Range("K22").Value = 855.7222
$ws.Range("L22").Value = 1456.4445
$ws.Range("M22").Value = -560.7222
$ws.Range("N22").Value = -2046.4445

$ws.Range("H27").Value = 1055.963
$ws.Range("I27").Value = 855.7222
$ws.Range("J27").Value = 1456.4445
$ws.Range("K27").Value = 855.7222
$ws.Range("L27").Value = 1456.4445
$ws.Range("M27").Value = -748.7222
$ws.Range("N27").Value = -1670.4445

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 55000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 55000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 55000
$ws.Range("N10").Value = -55338
$ws.Range("M10").ClearContents()
